$wb = $excel.ActiveWorkbook

$pwSheet = $wb.Worksheets.Item("PasswordRegression")
$homeSheet = $wb.Worksheets.Item("HomePageRegression")
$pwSheet.Move($homeSheet)

$acctSheet = $wb.Worksheets.Item("AccountSetup")
$acctSheet.Activate()
